# Fix formatting introduced when scraping floating point numbers.
#
# 1) Column H ("Importe") values were scraped in an ES/AR locale shape
#    ("." thousands separator, "," decimal separator), e.g. "44.395,00".
#    Re-emit them in plain decimal-point shape: "44395.00".
#
# 2) A handful of "Razon social" entries (column E) used a comma to
#    separate co-contracted people's names, which collided with the same
#    locale-formatting bug in the scraper and got mangled; replace the
#    stray commas with periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 295 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 8)
    $orig = $cell.Value2
    if ($orig -eq $null -or $orig -eq "") { continue }

    # Drop the "." thousands separators, then turn the "," decimal
    # separator into a ".".
    $fixed = $orig.Replace(".", "").Replace(",", ".")

    # Write back through NumberFormat "@" (Text) so Excel doesn't
    # re-interpret the digits as a real number (which would drop the
    # trailing zeros / thousands grouping again); ClearFormats()
    # afterwards removes the temporary Text format so the cell keeps its
    # original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value2 = $fixed
    $cell.ClearFormats()
}

# Razon social fixes: commas used as a name separator got corrupted by
# the same scraping bug -> replace with periods.
$ws.Range("E56").Value2  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E98").Value2  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E193").Value2 = "RICCOTTI. MARIANA EDITH"
$ws.Range("E199").Value2 = "DODERA. JORGE ABELARDO"
$ws.Range("E204").Value2 = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E221").Value2 = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
